# ====================================================================
# ZBP_13_dusevni_zdravi.xlsx -- add new wave "22. 2. 2022" of results,
# rename old "25. 1. 2021" header to "25. 1. 2022", refresh titles.
# ====================================================================

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "data" (percentages)
# ---------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("data")

# Fix mislabeled header: last column was "25. 1. 2021", should say 2022
$ws1.Range("AM1").Value = "25. 1. 2022"

# Append the new survey wave column (AN) with the same header style as the rest
$ws1.Range("AM1").Copy()
$ws1.Range("AN1").PasteSpecial(-4122)
$ws1.Range("AN1").Value = "22. 2. 2022"

# A few values in the "25. 1. 2022" column were revised with this update
$ws1.Range("AM6").Value = 0.13
$ws1.Range("AM7").Value = 0.11
$ws1.Range("AM10").Value = 0.13
$ws1.Range("AM13").Value = 0.13
$ws1.Range("AM15").Value = 0.18
$ws1.Range("AM17").Value = 0.09
$ws1.Range("AM18").Value = 0.11
$ws1.Range("AM19").Value = 0.09
$ws1.Range("AM20").Value = 0.21
$ws1.Range("AM23").Value = 0.19
$ws1.Range("AM24").Value = 0.11

# Populate the new "22. 2. 2022" column
$ws1.Range("AN2").Value = 0.1
$ws1.Range("AN3").Value = 0.09
$ws1.Range("AN4").Value = 0.11
$ws1.Range("AN5").Value = 0.08
$ws1.Range("AN6").Value = 0.11
$ws1.Range("AN7").Value = 0.11
$ws1.Range("AN8").Value = 0.1
$ws1.Range("AN9").Value = 0.08
$ws1.Range("AN10").Value = 0.13
$ws1.Range("AN11").Value = 0.05
$ws1.Range("AN12").Value = 0.08
$ws1.Range("AN13").Value = 0.15
$ws1.Range("AN14").Value = 0.12
$ws1.Range("AN15").Value = 0.13
$ws1.Range("AN16").Value = 0.11
$ws1.Range("AN17").Value = 0.09
$ws1.Range("AN18").Value = 0.12
$ws1.Range("AN19").Value = 0.09
$ws1.Range("AN20").Value = 0.2
$ws1.Range("AN21").Value = 0.12
$ws1.Range("AN22").Value = 0.09
$ws1.Range("AN23").Value = 0.18
$ws1.Range("AN24").Value = 0.12
$ws1.Range("AN25").Value = 0.07
$ws1.Range("AN26").Value = 0.1
$ws1.Range("AN27").Value = 0.04
$ws1.Range("AN28").Value = 0.08
$ws1.Range("AN29").Value = 0.23

# Refresh the footer title date stamp
$ws1.Range("A30").Value = "Život během pandemie, Duševní zdraví, % respondentů celkově a ve skupinách, aktualizace 2. 3. 2022"

# ---------------------------------------------------------------
# Sheet "pocetR" (respondent counts)
# ---------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("pocetR")

# Fix mislabeled header: last column was "25. 1. 2021", should say 2022
$ws2.Range("AL1").Value = "25. 1. 2022"

# Append the new survey wave column (AM) with the same header style as the rest
$ws2.Range("AL1").Copy()
$ws2.Range("AM1").PasteSpecial(-4122)
$ws2.Range("AM1").Value = "22. 2. 2022"

# The counts for "25. 1. 2022" are revised; the old values move to the new column
$ws2.Range("AL2").Value = 1848
$ws2.Range("AL3").Value = 419
$ws2.Range("AL4").Value = 1429
$ws2.Range("AL5").Value = 305
$ws2.Range("AL6").Value = 815
$ws2.Range("AL7").Value = 114
$ws2.Range("AL8").Value = 614
$ws2.Range("AL9").Value = 893
$ws2.Range("AL10").Value = 955
$ws2.Range("AL11").Value = 225
$ws2.Range("AL12").Value = 668
$ws2.Range("AL13").Value = 305
$ws2.Range("AL14").Value = 650
$ws2.Range("AL15").Value = 164
$ws2.Range("AL16").Value = 287
$ws2.Range("AL17").Value = 361
$ws2.Range("AL18").Value = 328
$ws2.Range("AL19").Value = 708
$ws2.Range("AL20").Value = 190
$ws2.Range("AL21").Value = 334
$ws2.Range("AL22").Value = 1324
$ws2.Range("AL23").Value = 172
$ws2.Range("AL24").Value = 642
$ws2.Range("AL25").Value = 619
$ws2.Range("AL26").Value = 307
$ws2.Range("AL27").Value = 499
$ws2.Range("AL28").Value = 773
$ws2.Range("AL29").Value = 576

# Populate the new "22. 2. 2022" column
$ws2.Range("AM2").Value = 1786
$ws2.Range("AM3").Value = 366
$ws2.Range("AM4").Value = 1420
$ws2.Range("AM5").Value = 275
$ws2.Range("AM6").Value = 808
$ws2.Range("AM7").Value = 91
$ws2.Range("AM8").Value = 612
$ws2.Range("AM9").Value = 863
$ws2.Range("AM10").Value = 923
$ws2.Range("AM11").Value = 207
$ws2.Range("AM12").Value = 656
$ws2.Range("AM13").Value = 293
$ws2.Range("AM14").Value = 630
$ws2.Range("AM15").Value = 159
$ws2.Range("AM16").Value = 286
$ws2.Range("AM17").Value = 346
$ws2.Range("AM18").Value = 317
$ws2.Range("AM19").Value = 678
$ws2.Range("AM20").Value = 170
$ws2.Range("AM21").Value = 348
$ws2.Range("AM22").Value = 1268
$ws2.Range("AM23").Value = 173
$ws2.Range("AM24").Value = 588
$ws2.Range("AM25").Value = 623
$ws2.Range("AM26").Value = 297
$ws2.Range("AM27").Value = 567
$ws2.Range("AM28").Value = 778
$ws2.Range("AM29").Value = 441

# Refresh the footer title date stamp
$ws2.Range("A30").Value = "Život během pandemie, Duševní zdraví, velikost dotázaného souboru celkově a ve skupinách, aktualizace 2. 3. 2022"

